# Auto-generated edit.ps1
# Applies the 2024-03 gh-pages data refresh to the "展览" and "全部类型"
# sheets of 江西-漫展信息.xlsx: a new row is inserted at the top of the
# event listing (上饶·囧喵喵次元国风动漫游戏展), which pushes every
# existing event down by one row (17 -> 18 data rows total), and several
# "想去人数" (interest count) values are refreshed to their latest scrape.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('展览')

# -- Column A: sequential row numbers 1..17, bordered/bold style --
for ($r = 2; $r -le 18; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}
# Row 18 is brand new -- copy A-column formatting (border/bold/center)
# from an existing data cell (A17 carries the same style as A1..A17).
$ws.Cells.Item(17, 1).Copy() | Out-Null
$ws.Cells.Item(18, 1).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# -- Columns B and E hold plain date-like text; force Text format so --
# -- Excel doesn't silently reinterpret them as date serials.        --
$dateLikeRange = $ws.Range('B2:B18')
$dateLikeRange.NumberFormat = '@'
$rangeRow = $ws.Range('E2:E18')
$rangeRow.NumberFormat = '@'

# Row 2
$ws.Cells.Item(2, 2).Value = '2024.02.23'
$ws.Cells.Item(2, 3).Value = '上饶·囧喵喵次元国风动漫游戏展'
$ws.Cells.Item(2, 4).Value = '春江北大道19号 博悦宴会艺术中心'
$ws.Cells.Item(2, 5).Value = '2024.02.23 09:00-02.23 17:00'
$ws.Cells.Item(2, 6).Value = 1093
$ws.Cells.Item(2, 7).Value = '不可售'
$ws.Cells.Item(2, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=80240'
$ws.Cells.Item(2, 9).Value = '//i0.hdslb.com/bfs/openplatform/202312/Qwh83wl31703836740097.jpeg'

# Row 3
$ws.Cells.Item(3, 2).Value = '2024.02.23'
$ws.Cells.Item(3, 3).Value = '南昌·国乙only·突破次元计划（取消）'
$ws.Cells.Item(3, 4).Value = '高处见美好生活公园 百家喜宴高新店'
$ws.Cells.Item(3, 5).Value = '2024.02.23 10:00-02.23 21:00'
$ws.Cells.Item(3, 6).Value = 306
$ws.Cells.Item(3, 7).Value = '不可售'
$ws.Cells.Item(3, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=80413'
$ws.Cells.Item(3, 9).Value = '//i0.hdslb.com/bfs/openplatform/202401/XvmB77wb1704252353395.jpeg'

# Row 4
$ws.Cells.Item(4, 2).Value = '2024.02.24'
$ws.Cells.Item(4, 3).Value = '南昌·Cookie动漫嘉年华-赵路专场票'
$ws.Cells.Item(4, 4).Value = '九龙大道1177号 南昌绿地国际博览中心'
$ws.Cells.Item(4, 5).Value = '2024.02.24 11:00-02.24 17:00'
$ws.Cells.Item(4, 6).Value = 363
$ws.Cells.Item(4, 7).Value = '已售罄'
$ws.Cells.Item(4, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=81769'
$ws.Cells.Item(4, 9).Value = '//i2.hdslb.com/bfs/openplatform/202402/DhCi2kWe1707123386859.png'

# Row 5
$ws.Cells.Item(5, 2).Value = '2024.02.24'
$ws.Cells.Item(5, 3).Value = '南昌·第一届Cookie动漫嘉年华'
$ws.Cells.Item(5, 4).Value = '九龙大道1177号 南昌绿地国际博览中心'
$ws.Cells.Item(5, 5).Value = '2024.02.24 09:00-02.24 17:00'
$ws.Cells.Item(5, 6).Value = 4686
$ws.Cells.Item(5, 7).Value = 65
$ws.Cells.Item(5, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=81033'
$ws.Cells.Item(5, 9).Value = '//i1.hdslb.com/bfs/openplatform/202401/P994oBkz1705562167665.png'

# Row 6
$ws.Cells.Item(6, 2).Value = '2024.02.24'
$ws.Cells.Item(6, 3).Value = '宜春·融荟城难忘今宵汉文化节'
$ws.Cells.Item(6, 4).Value = '宜阳大道239号 宜春融荟城'
$ws.Cells.Item(6, 5).Value = '2024.02.24 14:00-02.24 18:00'
$ws.Cells.Item(6, 6).Value = 27
$ws.Cells.Item(6, 7).Value = 10
$ws.Cells.Item(6, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=81690'
$ws.Cells.Item(6, 9).Value = '//i0.hdslb.com/bfs/openplatform/202402/ldtkc9Sp1706865634128.jpeg'

# Row 7
$ws.Cells.Item(7, 2).Value = '2024.02.24'
$ws.Cells.Item(7, 3).Value = '景德镇·陶溪川×次元文化元宵游园会（ 免费活动）'
$ws.Cells.Item(7, 4).Value = '新厂西路315号 陶溪川发布大厅'
$ws.Cells.Item(7, 5).Value = '2024.02.24 10:00-02.25 18:00'
$ws.Cells.Item(7, 6).Value = 402
$ws.Cells.Item(7, 7).Value = 30
$ws.Cells.Item(7, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=81207'
$ws.Cells.Item(7, 9).Value = '//i1.hdslb.com/bfs/openplatform/202402/nIs2jtUn1707298876430.png'

# Row 8
$ws.Cells.Item(8, 2).Value = '2024.03.02'
$ws.Cells.Item(8, 3).Value = '南昌·meeting动漫游戏嘉年华'
$ws.Cells.Item(8, 4).Value = '南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆'
$ws.Cells.Item(8, 5).Value = '2024.03.02 09:00-03.03 17:00'
$ws.Cells.Item(8, 6).Value = 1407
$ws.Cells.Item(8, 7).Value = 60
$ws.Cells.Item(8, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=79555'
$ws.Cells.Item(8, 9).Value = '//i0.hdslb.com/bfs/openplatform/202402/l6GUtggC1706843695971.jpeg'

# Row 9
$ws.Cells.Item(9, 2).Value = '2024.03.09'
$ws.Cells.Item(9, 3).Value = '景德镇·江报国风动漫展 '
$ws.Cells.Item(9, 4).Value = '迎宾大道与寺山路交叉口东200米 陶博城'
$ws.Cells.Item(9, 5).Value = '2024.03.09 09:00-03.10 17:00'
$ws.Cells.Item(9, 6).Value = 926
$ws.Cells.Item(9, 7).Value = 55
$ws.Cells.Item(9, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=81362'
$ws.Cells.Item(9, 9).Value = '//i2.hdslb.com/bfs/openplatform/202402/oM49o66R1708334630235.jpeg'

# Row 10
$ws.Cells.Item(10, 2).Value = '2024.03.16'
$ws.Cells.Item(10, 3).Value = '景德镇·原神X崩铁X崩坏动漫展only'
$ws.Cells.Item(10, 4).Value = '陶阳南路188号 晨枫臻品酒店'
$ws.Cells.Item(10, 5).Value = '2024.03.16 10:00-03.16 17:00'
$ws.Cells.Item(10, 6).Value = 56
$ws.Cells.Item(10, 7).Value = 55
$ws.Cells.Item(10, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=80920'
$ws.Cells.Item(10, 9).Value = '//i0.hdslb.com/bfs/openplatform/202401/IugBckTp1705469476482.png'

# Row 11
$ws.Cells.Item(11, 2).Value = '2024.03.16'
$ws.Cells.Item(11, 3).Value = '江西·ShiningStaR动漫游戏文化节5th'
$ws.Cells.Item(11, 4).Value = '高新开发区紫阳大道666号 江西奥林匹克体育中心综合训练馆'
$ws.Cells.Item(11, 5).Value = '2024.03.16 09:30-03.17 17:00'
$ws.Cells.Item(11, 6).Value = 1204
$ws.Cells.Item(11, 7).Value = 60
$ws.Cells.Item(11, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=81792'
$ws.Cells.Item(11, 9).Value = '//i2.hdslb.com/bfs/openplatform/202402/2l16aHBJ1707209383729.jpeg'

# Row 12
$ws.Cells.Item(12, 2).Value = '2024.03.23'
$ws.Cells.Item(12, 3).Value = '上饶·原×铁×崩only'
$ws.Cells.Item(12, 4).Value = '五三东大道42号 回禾酒店'
$ws.Cells.Item(12, 5).Value = '2024.03.23 10:00-03.23 17:00'
$ws.Cells.Item(12, 6).Value = 30
$ws.Cells.Item(12, 7).Value = 60
$ws.Cells.Item(12, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=81103'
$ws.Cells.Item(12, 9).Value = '//i2.hdslb.com/bfs/openplatform/202401/pp6c5TsC1705647180602.jpeg'

# Row 13
$ws.Cells.Item(13, 2).Value = '2024.03.23'
$ws.Cells.Item(13, 3).Value = '南昌·AP动漫游戏嘉年华'
$ws.Cells.Item(13, 4).Value = '八一桥街道青山南路118号 蓝海会展中心'
$ws.Cells.Item(13, 5).Value = '2024.03.23 09:00-03.24 17:00'
$ws.Cells.Item(13, 6).Value = 667
$ws.Cells.Item(13, 7).Value = 58.5
$ws.Cells.Item(13, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=81232'
$ws.Cells.Item(13, 9).Value = '//i2.hdslb.com/bfs/openplatform/202401/NZv97SmS1705912230957.jpeg'

# Row 14
$ws.Cells.Item(14, 2).Value = '2024.03.23'
$ws.Cells.Item(14, 3).Value = '南昌·原X穹X崩only'
$ws.Cells.Item(14, 4).Value = '丰和北大道299号 新吉花园酒店'
$ws.Cells.Item(14, 5).Value = '2024.03.23 10:00-03.23 17:00'
$ws.Cells.Item(14, 6).Value = 60
$ws.Cells.Item(14, 7).Value = 65
$ws.Cells.Item(14, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=80807'
$ws.Cells.Item(14, 9).Value = '//i0.hdslb.com/bfs/openplatform/202401/rY4v2Opx1705051458246.jpeg'

# Row 15
$ws.Cells.Item(15, 2).Value = '2024.03.23'
$ws.Cells.Item(15, 3).Value = '南昌·运动番only春季集训'
$ws.Cells.Item(15, 4).Value = '创新三路777号 南昌小飞侠章鱼文化体育公园'
$ws.Cells.Item(15, 5).Value = '2024.03.23 10:00-03.24 17:00'
$ws.Cells.Item(15, 6).Value = 54
$ws.Cells.Item(15, 7).Value = 58
$ws.Cells.Item(15, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=81950'
$ws.Cells.Item(15, 9).Value = '//i1.hdslb.com/bfs/openplatform/202402/bm4uH4qB1708425538357.jpeg'

# Row 16
$ws.Cells.Item(16, 2).Value = '2024.03.24'
$ws.Cells.Item(16, 3).Value = '南昌·AP动漫游戏  嘉年华内场票-小N&子音'
$ws.Cells.Item(16, 4).Value = '八一桥街道青山南路118号 蓝海会展中心'
$ws.Cells.Item(16, 5).Value = '2024.03.24 09:00-03.24 17:00'
$ws.Cells.Item(16, 6).Value = 27
$ws.Cells.Item(16, 7).Value = 218
$ws.Cells.Item(16, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=81973'
$ws.Cells.Item(16, 9).Value = '//i0.hdslb.com/bfs/openplatform/202402/zbG5HICL1708504962467.jpeg'

# Row 17
$ws.Cells.Item(17, 2).Value = '2024.03.30'
$ws.Cells.Item(17, 3).Value = '南昌·CM01动漫游戏博览会'
$ws.Cells.Item(17, 4).Value = '怀玉山大道1315号 南昌绿地国际博览中心'
$ws.Cells.Item(17, 5).Value = '2024.03.30 10:00-03.31 17:00'
$ws.Cells.Item(17, 6).Value = 284
$ws.Cells.Item(17, 7).Value = 55
$ws.Cells.Item(17, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=81691'
$ws.Cells.Item(17, 9).Value = '//i2.hdslb.com/bfs/openplatform/202402/IYLaH7AS1706866218597.png'

# Row 18
$ws.Cells.Item(18, 2).Value = '2024.03.30'
$ws.Cells.Item(18, 3).Value = '鹰潭·原×铁×崩only'
$ws.Cells.Item(18, 4).Value = '南站路24号 回禾酒店(鹰潭火车站店)'
$ws.Cells.Item(18, 5).Value = '2024.03.30 10:00-03.30 17:00'
$ws.Cells.Item(18, 6).Value = 24
$ws.Cells.Item(18, 7).Value = 60
$ws.Cells.Item(18, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=81097'
$ws.Cells.Item(18, 9).Value = '//i2.hdslb.com/bfs/openplatform/202401/q0AZaXAk1705646244207.jpeg'


$ws = $wb.Worksheets.Item('全部类型')

# -- Column A: sequential row numbers 1..17, bordered/bold style --
for ($r = 2; $r -le 18; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}
# Row 18 is brand new -- copy A-column formatting (border/bold/center)
# from an existing data cell (A17 carries the same style as A1..A17).
$ws.Cells.Item(17, 1).Copy() | Out-Null
$ws.Cells.Item(18, 1).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# -- Columns B and E hold plain date-like text; force Text format so --
# -- Excel doesn't silently reinterpret them as date serials.        --
$dateLikeRange = $ws.Range('B2:B18')
$dateLikeRange.NumberFormat = '@'
$rangeRow = $ws.Range('E2:E18')
$rangeRow.NumberFormat = '@'

# Row 2
$ws.Cells.Item(2, 2).Value = '2024.02.23'
$ws.Cells.Item(2, 3).Value = '上饶·囧喵喵次元国风动漫游戏展'
$ws.Cells.Item(2, 4).Value = '春江北大道19号 博悦宴会艺术中心'
$ws.Cells.Item(2, 5).Value = '2024.02.23 09:00-02.23 17:00'
$ws.Cells.Item(2, 6).Value = 1093
$ws.Cells.Item(2, 7).Value = '不可售'
$ws.Cells.Item(2, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=80240'
$ws.Cells.Item(2, 9).Value = '//i0.hdslb.com/bfs/openplatform/202312/Qwh83wl31703836740097.jpeg'

# Row 3
$ws.Cells.Item(3, 2).Value = '2024.02.23'
$ws.Cells.Item(3, 3).Value = '南昌·国乙only·突破次元计划（取消）'
$ws.Cells.Item(3, 4).Value = '高处见美好生活公园 百家喜宴高新店'
$ws.Cells.Item(3, 5).Value = '2024.02.23 10:00-02.23 21:00'
$ws.Cells.Item(3, 6).Value = 306
$ws.Cells.Item(3, 7).Value = '不可售'
$ws.Cells.Item(3, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=80413'
$ws.Cells.Item(3, 9).Value = '//i0.hdslb.com/bfs/openplatform/202401/XvmB77wb1704252353395.jpeg'

# Row 4
$ws.Cells.Item(4, 2).Value = '2024.02.24'
$ws.Cells.Item(4, 3).Value = '南昌·Cookie动漫嘉年华-赵路专场票'
$ws.Cells.Item(4, 4).Value = '九龙大道1177号 南昌绿地国际博览中心'
$ws.Cells.Item(4, 5).Value = '2024.02.24 11:00-02.24 17:00'
$ws.Cells.Item(4, 6).Value = 363
$ws.Cells.Item(4, 7).Value = '已售罄'
$ws.Cells.Item(4, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=81769'
$ws.Cells.Item(4, 9).Value = '//i2.hdslb.com/bfs/openplatform/202402/DhCi2kWe1707123386859.png'

# Row 5
$ws.Cells.Item(5, 2).Value = '2024.02.24'
$ws.Cells.Item(5, 3).Value = '南昌·第一届Cookie动漫嘉年华'
$ws.Cells.Item(5, 4).Value = '九龙大道1177号 南昌绿地国际博览中心'
$ws.Cells.Item(5, 5).Value = '2024.02.24 09:00-02.24 17:00'
$ws.Cells.Item(5, 6).Value = 4686
$ws.Cells.Item(5, 7).Value = 65
$ws.Cells.Item(5, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=81033'
$ws.Cells.Item(5, 9).Value = '//i1.hdslb.com/bfs/openplatform/202401/P994oBkz1705562167665.png'

# Row 6
$ws.Cells.Item(6, 2).Value = '2024.02.24'
$ws.Cells.Item(6, 3).Value = '宜春·融荟城难忘今宵汉文化节'
$ws.Cells.Item(6, 4).Value = '宜阳大道239号 宜春融荟城'
$ws.Cells.Item(6, 5).Value = '2024.02.24 14:00-02.24 18:00'
$ws.Cells.Item(6, 6).Value = 27
$ws.Cells.Item(6, 7).Value = 10
$ws.Cells.Item(6, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=81690'
$ws.Cells.Item(6, 9).Value = '//i0.hdslb.com/bfs/openplatform/202402/ldtkc9Sp1706865634128.jpeg'

# Row 7
$ws.Cells.Item(7, 2).Value = '2024.02.24'
$ws.Cells.Item(7, 3).Value = '景德镇·陶溪川×次元文化元宵游园会（ 免费活动）'
$ws.Cells.Item(7, 4).Value = '新厂西路315号 陶溪川发布大厅'
$ws.Cells.Item(7, 5).Value = '2024.02.24 10:00-02.25 18:00'
$ws.Cells.Item(7, 6).Value = 402
$ws.Cells.Item(7, 7).Value = 30
$ws.Cells.Item(7, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=81207'
$ws.Cells.Item(7, 9).Value = '//i1.hdslb.com/bfs/openplatform/202402/nIs2jtUn1707298876430.png'

# Row 8
$ws.Cells.Item(8, 2).Value = '2024.03.02'
$ws.Cells.Item(8, 3).Value = '南昌·meeting动漫游戏嘉年华'
$ws.Cells.Item(8, 4).Value = '南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆'
$ws.Cells.Item(8, 5).Value = '2024.03.02 09:00-03.03 17:00'
$ws.Cells.Item(8, 6).Value = 1407
$ws.Cells.Item(8, 7).Value = 60
$ws.Cells.Item(8, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=79555'
$ws.Cells.Item(8, 9).Value = '//i0.hdslb.com/bfs/openplatform/202402/l6GUtggC1706843695971.jpeg'

# Row 9
$ws.Cells.Item(9, 2).Value = '2024.03.09'
$ws.Cells.Item(9, 3).Value = '景德镇·江报国风动漫展 '
$ws.Cells.Item(9, 4).Value = '迎宾大道与寺山路交叉口东200米 陶博城'
$ws.Cells.Item(9, 5).Value = '2024.03.09 09:00-03.10 17:00'
$ws.Cells.Item(9, 6).Value = 926
$ws.Cells.Item(9, 7).Value = 55
$ws.Cells.Item(9, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=81362'
$ws.Cells.Item(9, 9).Value = '//i2.hdslb.com/bfs/openplatform/202402/oM49o66R1708334630235.jpeg'

# Row 10
$ws.Cells.Item(10, 2).Value = '2024.03.16'
$ws.Cells.Item(10, 3).Value = '景德镇·原神X崩铁X崩坏动漫展only'
$ws.Cells.Item(10, 4).Value = '陶阳南路188号 晨枫臻品酒店'
$ws.Cells.Item(10, 5).Value = '2024.03.16 10:00-03.16 17:00'
$ws.Cells.Item(10, 6).Value = 56
$ws.Cells.Item(10, 7).Value = 55
$ws.Cells.Item(10, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=80920'
$ws.Cells.Item(10, 9).Value = '//i0.hdslb.com/bfs/openplatform/202401/IugBckTp1705469476482.png'

# Row 11
$ws.Cells.Item(11, 2).Value = '2024.03.16'
$ws.Cells.Item(11, 3).Value = '江西·ShiningStaR动漫游戏文化节5th'
$ws.Cells.Item(11, 4).Value = '高新开发区紫阳大道666号 江西奥林匹克体育中心综合训练馆'
$ws.Cells.Item(11, 5).Value = '2024.03.16 09:30-03.17 17:00'
$ws.Cells.Item(11, 6).Value = 1205
$ws.Cells.Item(11, 7).Value = 60
$ws.Cells.Item(11, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=81792'
$ws.Cells.Item(11, 9).Value = '//i2.hdslb.com/bfs/openplatform/202402/2l16aHBJ1707209383729.jpeg'

# Row 12
$ws.Cells.Item(12, 2).Value = '2024.03.23'
$ws.Cells.Item(12, 3).Value = '上饶·原×铁×崩only'
$ws.Cells.Item(12, 4).Value = '五三东大道42号 回禾酒店'
$ws.Cells.Item(12, 5).Value = '2024.03.23 10:00-03.23 17:00'
$ws.Cells.Item(12, 6).Value = 30
$ws.Cells.Item(12, 7).Value = 60
$ws.Cells.Item(12, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=81103'
$ws.Cells.Item(12, 9).Value = '//i2.hdslb.com/bfs/openplatform/202401/pp6c5TsC1705647180602.jpeg'

# Row 13
$ws.Cells.Item(13, 2).Value = '2024.03.23'
$ws.Cells.Item(13, 3).Value = '南昌·AP动漫游戏嘉年华'
$ws.Cells.Item(13, 4).Value = '八一桥街道青山南路118号 蓝海会展中心'
$ws.Cells.Item(13, 5).Value = '2024.03.23 09:00-03.24 17:00'
$ws.Cells.Item(13, 6).Value = 667
$ws.Cells.Item(13, 7).Value = 58.5
$ws.Cells.Item(13, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=81232'
$ws.Cells.Item(13, 9).Value = '//i2.hdslb.com/bfs/openplatform/202401/NZv97SmS1705912230957.jpeg'

# Row 14
$ws.Cells.Item(14, 2).Value = '2024.03.23'
$ws.Cells.Item(14, 3).Value = '南昌·原X穹X崩only'
$ws.Cells.Item(14, 4).Value = '丰和北大道299号 新吉花园酒店'
$ws.Cells.Item(14, 5).Value = '2024.03.23 10:00-03.23 17:00'
$ws.Cells.Item(14, 6).Value = 60
$ws.Cells.Item(14, 7).Value = 65
$ws.Cells.Item(14, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=80807'
$ws.Cells.Item(14, 9).Value = '//i0.hdslb.com/bfs/openplatform/202401/rY4v2Opx1705051458246.jpeg'

# Row 15
$ws.Cells.Item(15, 2).Value = '2024.03.23'
$ws.Cells.Item(15, 3).Value = '南昌·运动番only春季集训'
$ws.Cells.Item(15, 4).Value = '创新三路777号 南昌小飞侠章鱼文化体育公园'
$ws.Cells.Item(15, 5).Value = '2024.03.23 10:00-03.24 17:00'
$ws.Cells.Item(15, 6).Value = 54
$ws.Cells.Item(15, 7).Value = 58
$ws.Cells.Item(15, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=81950'
$ws.Cells.Item(15, 9).Value = '//i1.hdslb.com/bfs/openplatform/202402/bm4uH4qB1708425538357.jpeg'

# Row 16
$ws.Cells.Item(16, 2).Value = '2024.03.24'
$ws.Cells.Item(16, 3).Value = '南昌·AP动漫游戏  嘉年华内场票-小N&子音'
$ws.Cells.Item(16, 4).Value = '八一桥街道青山南路118号 蓝海会展中心'
$ws.Cells.Item(16, 5).Value = '2024.03.24 09:00-03.24 17:00'
$ws.Cells.Item(16, 6).Value = 27
$ws.Cells.Item(16, 7).Value = 218
$ws.Cells.Item(16, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=81973'
$ws.Cells.Item(16, 9).Value = '//i0.hdslb.com/bfs/openplatform/202402/zbG5HICL1708504962467.jpeg'

# Row 17
$ws.Cells.Item(17, 2).Value = '2024.03.30'
$ws.Cells.Item(17, 3).Value = '南昌·CM01动漫游戏博览会'
$ws.Cells.Item(17, 4).Value = '怀玉山大道1315号 南昌绿地国际博览中心'
$ws.Cells.Item(17, 5).Value = '2024.03.30 10:00-03.31 17:00'
$ws.Cells.Item(17, 6).Value = 284
$ws.Cells.Item(17, 7).Value = 55
$ws.Cells.Item(17, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=81691'
$ws.Cells.Item(17, 9).Value = '//i2.hdslb.com/bfs/openplatform/202402/IYLaH7AS1706866218597.png'

# Row 18
$ws.Cells.Item(18, 2).Value = '2024.03.30'
$ws.Cells.Item(18, 3).Value = '鹰潭·原×铁×崩only'
$ws.Cells.Item(18, 4).Value = '南站路24号 回禾酒店(鹰潭火车站店)'
$ws.Cells.Item(18, 5).Value = '2024.03.30 10:00-03.30 17:00'
$ws.Cells.Item(18, 6).Value = 24
$ws.Cells.Item(18, 7).Value = 60
$ws.Cells.Item(18, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=81097'
$ws.Cells.Item(18, 9).Value = '//i2.hdslb.com/bfs/openplatform/202401/q0AZaXAk1705646244207.jpeg'

